$d = $word.ActiveDocument

# Step 1: Replace "ClosedDoer" with "ClosedDoor"
$d.Content.Find.Execute("ClosedDoer", $true, $false, $false, $false, $false,
                         $true, 1, $false, "ClosedDoor", 2)
